# Modified Report and Test Cases
# Adds a new "Sheet2" worksheet with extended test data (roles, profile,
# admin-content, messages) and tweaks the existing Sheet1 selection.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Add Sheet2 right after Sheet1
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# ---------------------------------------------------------------------
# 2. Reuse Sheet1's existing header / label / hyperlink styles for the
#    columns that mirror Sheet1 (A = username, B = password)
# ---------------------------------------------------------------------
$ws1.Range("A1:B1").Copy()
$ws2.Range("A1:B1").PasteSpecial(-4122)

$ws1.Range("B2").Copy()
$ws2.Range("B2").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 3. New "sub header" style (row 1, cols C..L): yellow fill + thin border
# ---------------------------------------------------------------------
$hdr = $ws2.Range("C1:L1")
$hdr.Interior.Color = 65535
$hdr.Borders.Color = 0
$hdr.Borders.LineStyle = 1

# ---------------------------------------------------------------------
# 4. New "data" style (row 2): thin border, no fill
# ---------------------------------------------------------------------
$data1 = $ws2.Range("C2:H2")
$data1.Borders.Color = 0
$data1.Borders.LineStyle = 1

$data2 = $ws2.Range("J2:K2")
$data2.Borders.Color = 0
$data2.Borders.LineStyle = 1

$data3 = $ws2.Range("L2")
$data3.Borders.Color = 0
$data3.Borders.LineStyle = 1

# Date-of-birth cell: same border, plus a date number format
$dob = $ws2.Range("I2")
$dob.Borders.Color = 0
$dob.Borders.LineStyle = 1
$dob.NumberFormat = "mm-dd-yy"

# ---------------------------------------------------------------------
# 5. Cell values -- written in the exact order the new strings were
#    first introduced so the shared-string table lines up.
# ---------------------------------------------------------------------
$ws2.Range("A1").Value = "username"
$ws2.Range("B1").Value = "password"
$ws2.Range("C1").Value = "pwd1"
$ws2.Range("C2").Value = "test123"
$ws2.Range("D1").Value = "CustomerRoles"
$ws2.Range("D2").Value = "Guest"
$ws2.Range("E1").Value = "ManagerOfVendor"
$ws2.Range("E2").Value = "Vendor 2"
$ws2.Range("F1").Value = "Gender"
$ws2.Range("F2").Value = "Male"
$ws2.Range("G1").Value = "FirstName"
$ws2.Range("H1").Value = "LastName"
$ws2.Range("G2").Value = "Jagadeesh"
$ws2.Range("H2").Value = "Patil"
$ws2.Range("I1").Value = "Dob"
$ws2.Range("J1").Value = "CompanyName"
$ws2.Range("J2").Value = "busyQA"
$ws2.Range("K1").Value = "AdminContent"
$ws2.Range("K2").Value = "This is for testing........."
$ws2.Range("L2").Value = "The new customer has been added successfully"
$ws2.Range("L1").Value = "Exmessage"

$ws2.Range("A2").Value = "admin@yourstore.com "
$ws2.Range("B2").Value = "admin"
$ws2.Range("I2").Value = [DateTime]"1985-07-05"

# ---------------------------------------------------------------------
# 6. Hyperlink for the admin e-mail in column A
# ---------------------------------------------------------------------
$ws2.Hyperlinks.Add($ws2.Range("A2"), "mailto:admin@yourstore.com") | Out-Null

# ---------------------------------------------------------------------
# 7. Column widths / row heights to roughly match the authored layout
# ---------------------------------------------------------------------
$ws2.Columns.Item(1).ColumnWidth = 20.7265625
$ws2.Columns.Item(2).ColumnWidth = 14.26953125
$ws2.Columns.Item(4).ColumnWidth = 13.453125
$ws2.Columns.Item(5).ColumnWidth = 19
$ws2.Columns.Item(10).ColumnWidth = 13.7265625
$ws2.Columns.Item(11).ColumnWidth = 19.90625
$ws2.Columns.Item(12).ColumnWidth = 40.7265625

$ws2.Range("A1:L2").RowHeight = 23.5

# ---------------------------------------------------------------------
# 8. Sheet1: move the selection from B3 to B2 (do this before switching
#    the active sheet, since selecting a range activates its sheet)
# ---------------------------------------------------------------------
$ws1.Range("B2").Select()

# ---------------------------------------------------------------------
# 9. Selection / activation -- Sheet2 becomes the active tab
# ---------------------------------------------------------------------
$ws2.Range("C1").Select()
$ws2.Activate()
